# The document's 3rd table is the cable-insulation-resistance table:
#   columns: [marking] [mark] [mark] [insulation resistance, MOhm] [note]
# Each data row (rows 2..13, since row 1 is the header) has its resistance
# value in column 3. Update those values per the target diff.

$d = $word.ActiveDocument

$tbl = $d.Tables.Item(3)

# 1-based table Row -> old/new resistance value (column 3)
$changes = @(
    @{ Row = 2;  Old = "0.094"; New = "0.059" }
    @{ Row = 3;  Old = "0.519"; New = "0.702" }
    @{ Row = 4;  Old = "0.924"; New = "0.438" }
    @{ Row = 5;  Old = "0.754"; New = "0.314" }
    @{ Row = 6;  Old = "0.351"; New = "0.606" }
    @{ Row = 7;  Old = "0.842"; New = "0.162" }
    @{ Row = 8;  Old = "0.068"; New = "0.192" }
    @{ Row = 9;  Old = "0.677"; New = "0.52"  }
    @{ Row = 10; Old = "0.236"; New = "0.846" }
    @{ Row = 11; Old = "0.89";  New = "0.968" }
    @{ Row = 12; Old = "0.89";  New = "0.141" }
    @{ Row = 13; Old = "0.198"; New = "0.233" }
)

foreach ($chg in $changes) {
    $cell = $tbl.Cell($chg.Row, 3)
    $rng = $cell.Range
    # wdReplaceOne (the last arg, 1) so the substitution is confined to this
    # single cell's Range and does not spill into sibling cells that share
    # the same old value (e.g. the two "0.89" rows).
    $ok = $rng.Find.Execute($chg.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $chg.New, 1)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for row $($chg.Row) ($($chg.Old) -> $($chg.New))"
    }
}

Write-Output "done"
